# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): copy the existing header formatting (bold,
# centered, bordered) from AA1:AC1 onto the three new header cells, then
# set their labels.
$ws.Range("AA1:AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2-57): every player on the roster shares the team's 1996
# season record (85 wins, 77 losses, 0 ties).
$ws.Range("AD2:AD57").Value = 85
$ws.Range("AE2:AE57").Value = 77
$ws.Range("AF2:AF57").Value = 0
